$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1575.5
$ws.Range("I9").Value = 1984.3334
$ws.Range("K9").Value = 1984.3334
$ws.Range("M9").Value = -1815.3334

$ws.Range("H64").Value = 5495.6665
$ws.Range("J64").Value = 5495.6665
$ws.Range("L64").Value = 5495.6665
$ws.Range("N64").Value = -5991.6665

$ws.Range("H67").Value = 5495.6665
$ws.Range("J67").Value = 5495.6665
$ws.Range("L67").Value = 5495.6665
$ws.Range("N67").Value = -7211.6665

$ws.Range("H69").Value = 35857.855
$ws.Range("I69").Value = 118505
$ws.Range("K69").Value = 355515
$ws.Range("M69").Value = -354641

$ws.Range("H72").Value = 35857.855
$ws.Range("I72").Value = 118505
$ws.Range("K72").Value = 1066545
$ws.Range("M72").Value = -1062177

$ws.Range("H88").Value = 1999
$ws.Range("J88").Value = 1999
$ws.Range("L88").Value = 1999
$ws.Range("N88").Value = -2811

$ws.Range("H91").Value = 1999
$ws.Range("J91").Value = 1999
$ws.Range("L91").Value = 1999
$ws.Range("N91").Value = -4807

$ws.Range("H100").Value = 3666.3333
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459

$ws.Range("H137").Value = 1899.7778
$ws.Range("I137").Value = 1019.6
$ws.Range("K137").Value = 3058.8
$ws.Range("M137").Value = -508.8000000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 690
$ws.Range("I74").Value = 916.6667
$ws.Range("J74").Value = 350
$ws.Range("K74").Value = 916.6667
$ws.Range("L74").Value = 350
$ws.Range("M74").Value = -42.66669999999999
$ws.Range("N74").Value = -2098

$ws.Range("H77").Value = 690
$ws.Range("I77").Value = 916.6667
$ws.Range("J77").Value = 350
$ws.Range("K77").Value = 4583.3335
$ws.Range("L77").Value = 1750
$ws.Range("M77").Value = -215.3334999999997
$ws.Range("N77").Value = -10486

$ws.Range("H110").Value = 1617.5
$ws.Range("I110").Value = 1617.5
$ws.Range("K110").Value = 1617.5
$ws.Range("M110").Value = 427.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2420.6365
$ws.Range("I105").Value = 2304.1428
$ws.Range("J105").Value = 2624.5
$ws.Range("K105").Value = 2304.1428
$ws.Range("L105").Value = 2624.5
$ws.Range("M105").Value = -557.1428000000001
$ws.Range("N105").Value = -6118.5

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 2518.182
$ws.Range("I134").Value = 2518.182
$ws.Range("K134").Value = 7554.545999999999
$ws.Range("M134").Value = -5019.545999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1360.2
$ws.Range("I31").Value = 1230.2858
$ws.Range("J31").Value = 1663.3334
$ws.Range("K31").Value = 1230.2858
$ws.Range("L31").Value = 1663.3334
$ws.Range("M31").Value = -935.2858000000001
$ws.Range("N31").Value = -2253.3334

$ws.Range("H34").Value = 1360.2
$ws.Range("I34").Value = 1230.2858
$ws.Range("J34").Value = 1663.3334
$ws.Range("K34").Value = 1230.2858
$ws.Range("L34").Value = 1663.3334
$ws.Range("M34").Value = -1028.2858
$ws.Range("N34").Value = -2067.3334

$ws.Range("H70").Value = 20000
$ws.Range("J70").Value = 20000
$ws.Range("L70").Value = 20000
$ws.Range("N70").Value = -20630

$ws.Range("H73").Value = 20000
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 20000
$ws.Range("N73").Value = -22184

$ws.Range("H132").Value = 2861.3333
$ws.Range("J132").Value = 2070.25
$ws.Range("L132").Value = 6210.75
$ws.Range("N132").Value = -11270.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2450
$ws.Range("I109").Value = 687.5
$ws.Range("J109").Value = 4800
$ws.Range("K109").Value = 2062.5
$ws.Range("L109").Value = 14400
$ws.Range("M109").Value = -1022.5
$ws.Range("N109").Value = -16480

$ws.Range("H129").Value = 1253685.2
$ws.Range("J129").Value = 2004206.4
$ws.Range("L129").Value = 6012619.199999999
$ws.Range("N129").Value = -6022619.199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 101.125
$ws.Range("I2").Value = 118.57895
$ws.Range("J2").Value = 34.8
$ws.Range("K2").Value = 118.57895
$ws.Range("L2").Value = 34.8
$ws.Range("M2").Value = -5.578950000000006
$ws.Range("N2").Value = -260.8

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H99").Value = 12349
$ws.Range("I99").Value = 12624.5
$ws.Range("K99").Value = 12624.5
$ws.Range("M99").Value = -10378.5

$ws.Range("H132").Value = 1727.1818
$ws.Range("I132").Value = 1727.1818
$ws.Range("K132").Value = 5181.5454
$ws.Range("M132").Value = -2651.5454

$ws.Range("H134").Value = 53721
$ws.Range("J134").Value = 53721
$ws.Range("L134").Value = 161163
$ws.Range("N134").Value = -166233

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3319.2856
$ws.Range("I40").Value = 2447
$ws.Range("K40").Value = 2447
$ws.Range("M40").Value = -2311

$ws.Range("H68").Value = 2600
$ws.Range("I68").Value = 2600
$ws.Range("K68").Value = 2600
$ws.Range("M68").Value = -1851

$ws.Range("H71").Value = 2600
$ws.Range("I71").Value = 2600
$ws.Range("K71").Value = 13000
$ws.Range("M71").Value = -9256

$ws.Range("H93").Value = 3159
$ws.Range("I93").Value = 2883
$ws.Range("K93").Value = 2883
$ws.Range("M93").Value = -1635

$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H135").Value = 89999.5
$ws.Range("J135").Value = 89999.5
$ws.Range("L135").Value = 89999.5
$ws.Range("N135").Value = -100139.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 10000
$ws.Range("J24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("N24").Value = -10460

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H41").Value = 19981.6
$ws.Range("I41").Value = 19973.5
$ws.Range("J41").Value = 19987
$ws.Range("K41").Value = 19973.5
$ws.Range("L41").Value = 19987
$ws.Range("M41").Value = -19583.5
$ws.Range("N41").Value = -20767

$ws.Range("H75").Value = 73201.60000000001
$ws.Range("J75").Value = 73333.336
$ws.Range("L75").Value = 73333.336
$ws.Range("N75").Value = -75205.336

$ws.Range("H78").Value = 73201.60000000001
$ws.Range("J78").Value = 73333.336
$ws.Range("L78").Value = 220000.008
$ws.Range("N78").Value = -229360.008

$ws.Range("H132").Value = 1472.04
$ws.Range("I132").Value = 1480.7142
$ws.Range("J132").Value = 1426.5
$ws.Range("K132").Value = 4279.5
$ws.Range("M132").Value = -1912.142599999999
$ws.Range("N132").Value = -9339.5
